$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 9 (which holds 99999/99999),
# shifting it down to row 10, and fill the new row 9 with 99993/99993.
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = 99993
$ws.Range("B9").Value = 99993
